$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 70, pushing existing rows 70-111 down to 72-113.
$ws.Rows("70:71").Insert()

# --- Fill in new row 70 ---
$ws.Cells.Item(70, 1).Value  = 5
$ws.Cells.Item(70, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(70, 3).Value  = "Maule"
$ws.Cells.Item(70, 4).Value  = 44957
$ws.Cells.Item(70, 5).Value  = 7
$ws.Cells.Item(70, 6).Value  = "Fruta"
$ws.Cells.Item(70, 7).Value  = 100101
$ws.Cells.Item(70, 8).Value  = "Berries"
$ws.Cells.Item(70, 9).Value  = 100101001
$ws.Cells.Item(70, 10).Value = "Arándano (blue)"
$ws.Cells.Item(70, 11).Value = "Sin especificar"
$ws.Cells.Item(70, 12).Value = "Primera"
$ws.Cells.Item(70, 13).Value = 200
$ws.Cells.Item(70, 14).Value = 3000
$ws.Cells.Item(70, 15).Value = 3000
$ws.Cells.Item(70, 16).Value = 3000
$ws.Cells.Item(70, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(70, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(70, 19).Value = 1500
$ws.Cells.Item(70, 20).Value = 2

# --- Fill in new row 71 ---
$ws.Cells.Item(71, 1).Value  = 5
$ws.Cells.Item(71, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(71, 3).Value  = "Maule"
$ws.Cells.Item(71, 4).Value  = 44957
$ws.Cells.Item(71, 5).Value  = 7
$ws.Cells.Item(71, 6).Value  = "Fruta"
$ws.Cells.Item(71, 7).Value  = 100101
$ws.Cells.Item(71, 8).Value  = "Berries"
$ws.Cells.Item(71, 9).Value  = 100101001
$ws.Cells.Item(71, 10).Value = "Arándano (blue)"
$ws.Cells.Item(71, 11).Value = "Sin especificar"
$ws.Cells.Item(71, 12).Value = "Segunda"
$ws.Cells.Item(71, 13).Value = 200
$ws.Cells.Item(71, 14).Value = 2500
$ws.Cells.Item(71, 15).Value = 2500
$ws.Cells.Item(71, 16).Value = 2500
$ws.Cells.Item(71, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(71, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(71, 19).Value = 1250
$ws.Cells.Item(71, 20).Value = 2
